{"js": "// Remove the direct (paragraph-level) <w:contextualSpacing .../> element\n// from every paragraph's properties (w:pPr) in the document body.\n//\n// The Word JS API does not expose a dedicated \"contextualSpacing\" getter/\n// setter on Word.Paragraph / Word.ParagraphFormat in this host, so we do a\n// surgical per-paragraph OOXML round-trip: read each paragraph's OOXML,\n// strip the <w:contextualSpacing/> element, and write the paragraph's XML\n// back via its own Range (Range.insertOoxml supports \"Replace\", unlike\n// Paragraph.insertOoxml).\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\nconst contextualSpacingRe = /<w:contextualSpacing\\b[^>]*\\/>/g;\n\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  const paragraph = paragraphs.items[i];\n\n  const ooxml = paragraph.getOoxml();\n  await context.sync();\n\n  const xml = ooxml.value;\n  if (!xml || xml.indexOf(\"w:contextualSpacing\") === -1) {\n    continue;\n  }\n\n  const updatedXml = xml.replace(contextualSpacingRe, \"\");\n\n  const range = paragraph.getRange();\n  range.insertOoxml(updatedXml, Word.InsertLocation.replace);\n  await context.sync();\n}\n", "ps1": "# Remove the direct (paragraph-level) <w:contextualSpacing .../> element\n# from every paragraph's properties (w:pPr) in the document body.\n#\n# The Word object model does not expose a dedicated ContextualSpacing\n# property on Paragraph / ParagraphFormat in this host, so we do a\n# surgical per-paragraph WordOpenXML round-trip: read each paragraph's\n# Range.XML(), strip the <w:contextualSpacing/> element, and write it\n# back with Range.InsertXML() (which replaces just that range's content).\n\n$d = $word.ActiveDocument\n$paragraphs = $d.Paragraphs\n$count = $paragraphs.Count\n\nfor ($i = 1; $i -le $count; $i++) {\n    $paragraph = $paragraphs.Item($i)\n    $range = $paragraph.Range\n    $xml = $range.XML()\n\n    if ($xml -notlike \"*contextualSpacing*\") {\n        continue\n    }\n\n    $updatedXml = $xml -replace '<w:contextualSpacing[^>]*/>', ''\n    [void]$range.InsertXML($updatedXml)\n}\n"}
